# Auto-generated Excel COM-interop script to apply Zalera_Profits.xlsx diff
# Updates cached market-price / profit columns (H-N) across 8 sheets as per scheduled runner refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 7505.6665
$ws.Range("J45").Value = 6250
$ws.Range("L45").Value = 18750
$ws.Range("N45").Value = -19134

$ws.Range("H69").Value = 7330
$ws.Range("J69").Value = 10000
$ws.Range("L69").Value = 30000
$ws.Range("N69").Value = -31748

$ws.Range("H72").Value = 7330
$ws.Range("J72").Value = 10000
$ws.Range("L72").Value = 90000
$ws.Range("N72").Value = -98736

$ws.Range("H133").Value = 115000
$ws.Range("I133").Value = 110000
$ws.Range("J133").Value = 120000
$ws.Range("K133").Value = 110000
$ws.Range("L133").Value = 120000
$ws.Range("M133").Value = -104940
$ws.Range("N133").Value = -130120

$ws.Range("H141").Value = 2599
$ws.Range("I141").Value = 1998.75
$ws.Range("K141").Value = 5996.25
$ws.Range("M141").Value = -816.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5780.08
$ws.Range("I132").Value = 3768.4
$ws.Range("J132").Value = 8797.6
$ws.Range("K132").Value = 11305.2
$ws.Range("L132").Value = 26392.8
$ws.Range("M132").Value = -8775.200000000001
$ws.Range("N132").Value = -31452.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 27800.715
$ws.Range("I102").Value = 13498.75
$ws.Range("K102").Value = 13498.75
$ws.Range("M102").Value = -10253.75

$ws.Range("H107").Value = 2043.45
$ws.Range("I107").Value = 1885.4
$ws.Range("J107").Value = 2517.6
$ws.Range("K107").Value = 1885.4
$ws.Range("L107").Value = 2517.6
$ws.Range("M107").Value = 34.59999999999991
$ws.Range("N107").Value = -6357.6

$ws.Range("H132").Value = 119769.336
$ws.Range("J132").Value = 119769.336
$ws.Range("L132").Value = 119769.336
$ws.Range("N132").Value = -129889.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 291.83334
$ws.Range("I7").Value = 133.55556
$ws.Range("J7").Value = 766.6667
$ws.Range("K7").Value = 133.55556
$ws.Range("L7").Value = 766.6667
$ws.Range("M7").Value = -20.55556000000001
$ws.Range("N7").Value = -992.6667

$ws.Range("H28").Value = 30178.75
$ws.Range("J28").Value = 30178.75
$ws.Range("L28").Value = 30178.75
$ws.Range("N28").Value = -30668.75

$ws.Range("H31").Value = 41670908
$ws.Range("I31").Value = 333334270
$ws.Range("J31").Value = 4714.4287
$ws.Range("K31").Value = 333334270
$ws.Range("L31").Value = 4714.4287
$ws.Range("M31").Value = -333333975
$ws.Range("N31").Value = -5304.4287

$ws.Range("H34").Value = 41670908
$ws.Range("I34").Value = 333334270
$ws.Range("J34").Value = 4714.4287
$ws.Range("K34").Value = 333334270
$ws.Range("L34").Value = 4714.4287
$ws.Range("M34").Value = -333334068
$ws.Range("N34").Value = -5118.4287

$ws.Range("H62").Value = 9704.666999999999
$ws.Range("I62").Value = 4972.25
$ws.Range("J62").Value = 15113.143
$ws.Range("K62").Value = 4972.25
$ws.Range("L62").Value = 15113.143
$ws.Range("M62").Value = -4348.25
$ws.Range("N62").Value = -16361.143

$ws.Range("H65").Value = 9704.666999999999
$ws.Range("I65").Value = 4972.25
$ws.Range("J65").Value = 15113.143
$ws.Range("K65").Value = 24861.25
$ws.Range("L65").Value = 75565.715
$ws.Range("M65").Value = -21741.25
$ws.Range("N65").Value = -81805.715

$ws.Range("H107").Value = 1642
$ws.Range("I107").Value = 1397.7368
$ws.Range("K107").Value = 1397.7368
$ws.Range("M107").Value = 522.2632000000001

$ws.Range("H134").Value = 8555.385
$ws.Range("I134").Value = 8465.6
$ws.Range("K134").Value = 25396.8
$ws.Range("M134").Value = -22861.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 1935.4
$ws.Range("J69").Value = 3851.6667
$ws.Range("L69").Value = 11555.0001
$ws.Range("N69").Value = -13177.0001

$ws.Range("H72").Value = 1935.4
$ws.Range("J72").Value = 3851.6667
$ws.Range("L72").Value = 34665.0003
$ws.Range("N72").Value = -42777.0003

$ws.Range("H92").Value = 1380.0714
$ws.Range("I92").Value = 8002
$ws.Range("J92").Value = 870.6923
$ws.Range("K92").Value = 24006
$ws.Range("L92").Value = 2612.0769
$ws.Range("M92").Value = -22758
$ws.Range("N92").Value = -5108.0769

$ws.Range("H132").Value = 1998.1666
$ws.Range("I132").Value = 1998.1666
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 17983.4994
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -15453.4994
$ws.Range("N132").ClearContents()

$ws.Range("H134").Value = 250000500
$ws.Range("I134").Value = 250000500
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 750001500
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -749996430
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

$ws.Range("H137").Value = 10000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 10000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 30000
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -40200

$ws.Range("H139").Value = 6481.3335
$ws.Range("J139").Value = 6481.3335
$ws.Range("L139").Value = 19444.0005
$ws.Range("N139").Value = -29724.0005

$ws.Range("H140").Value = 41667868
$ws.Range("J140").Value = 1500
$ws.Range("L140").Value = 4500
$ws.Range("N140").Value = -14860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13191.818
$ws.Range("I70").Value = 12676.23
$ws.Range("K70").Value = 12676.23
$ws.Range("M70").Value = -12406.23

$ws.Range("H73").Value = 13191.818
$ws.Range("I73").Value = 12676.23
$ws.Range("K73").Value = 12676.23
$ws.Range("M73").Value = -11740.23

$ws.Range("H113").Value = 46873.25
$ws.Range("I113").Value = 4997
$ws.Range("J113").Value = 71999
$ws.Range("K113").Value = 4997
$ws.Range("L113").Value = 71999
$ws.Range("M113").Value = -2827
$ws.Range("N113").Value = -76339

$ws.Range("H132").Value = 6045.5
$ws.Range("I132").Value = 3663.45
$ws.Range("K132").Value = 10990.35
$ws.Range("M132").Value = -8460.349999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4440.5835
$ws.Range("I40").Value = 4664.4546
$ws.Range("K40").Value = 4664.4546
$ws.Range("M40").Value = -4528.4546

$ws.Range("H74").Value = 48786.75
$ws.Range("I74").Value = 48786.75
$ws.Range("K74").Value = 48786.75
$ws.Range("M74").Value = -47788.75

$ws.Range("H77").Value = 48786.75
$ws.Range("I77").Value = 48786.75
$ws.Range("K77").Value = 146360.25
$ws.Range("M77").Value = -141368.25

$ws.Range("H132").Value = 6494.7026
$ws.Range("I132").Value = 5985.4585
$ws.Range("K132").Value = 17956.3755
$ws.Range("M132").Value = -15426.3755

$ws.Range("H136").Value = 4369.241
$ws.Range("I136").Value = 3042.4285
$ws.Range("K136").Value = 9127.2855
$ws.Range("M136").Value = -6577.2855

$ws.Range("H139").Value = 120000
$ws.Range("J139").Value = 120000
$ws.Range("L139").Value = 120000
$ws.Range("N139").Value = -130280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 13452
$ws.Range("J68").Value = 13452
$ws.Range("L68").Value = 13452
$ws.Range("N68").Value = -15074

$ws.Range("H71").Value = 13452
$ws.Range("J71").Value = 13452
$ws.Range("L71").Value = 40356
$ws.Range("N71").Value = -48468

$ws.Range("H124").Value = 110067.5
$ws.Range("J124").Value = 110067.5
$ws.Range("L124").Value = 110067.5
$ws.Range("N124").Value = -119887.5

$ws.Range("H132").Value = 3529.9556
$ws.Range("I132").Value = 1156.5555
$ws.Range("J132").Value = 7090.0557
$ws.Range("K132").Value = 3469.6665
$ws.Range("L132").Value = 21270.1671
$ws.Range("M132").Value = -939.6664999999998
$ws.Range("N132").Value = -26330.1671

$ws.Range("H136").Value = 2932.5667
$ws.Range("I136").Value = 2118.65
$ws.Range("J136").Value = 4560.4
$ws.Range("K136").Value = 6355.950000000001
$ws.Range("L136").Value = 13681.2
$ws.Range("M136").Value = -3805.950000000001
$ws.Range("N136").Value = -18781.2
